# Update "想去人数" (F) / "最低票价" (G) figures across the four sheets to
# match the newly scraped numbers (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 195
$ws.Range("F6").Value = 839
$ws.Range("F7").Value = 4248
$ws.Range("F11").Value = 6211
$ws.Range("F12").Value = 6211
$ws.Range("F14").Value = 469
$ws.Range("F15").Value = 2381
$ws.Range("F17").Value = 171
$ws.Range("F18").Value = 490
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = 9346
$ws.Range("F22").Value = 2522
$ws.Range("F24").Value = 2336
$ws.Range("F25").Value = 2493
$ws.Range("F28").Value = 1987
$ws.Range("F31").Value = 340
$ws.Range("F33").Value = 52
$ws.Range("F34").Value = 283
$ws.Range("F35").Value = 47
$ws.Range("F36").Value = 92
$ws.Range("F38").Value = 1233
$ws.Range("F40").Value = 77
$ws.Range("F41").Value = 103
$ws.Range("F42").Value = 1572
$ws.Range("F43").Value = 2592
$ws.Range("F45").Value = 936
$ws.Range("F46").Value = 317
$ws.Range("F47").Value = 1257
$ws.Range("F48").Value = 30
$ws.Range("F50").Value = 9

# ---- 演出 (Performance) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G6").Value = "不可售"
$ws.Range("F10").Value = 933
$ws.Range("F12").Value = 153
$ws.Range("F22").Value = 113

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 701
$ws.Range("F3").Value = 916

# ---- 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 195
$ws.Range("F3").Value = 701
$ws.Range("F4").Value = 916
$ws.Range("F10").Value = 839
$ws.Range("F11").Value = 4248
$ws.Range("F16").Value = 6211
$ws.Range("F18").Value = 2381
$ws.Range("F19").Value = 171
$ws.Range("F20").Value = 490
$ws.Range("F21").Value = 9346
$ws.Range("F22").Value = 153
$ws.Range("F24").Value = 2522
$ws.Range("F26").Value = 2493
$ws.Range("F28").Value = 1987
$ws.Range("F31").Value = 340
$ws.Range("F33").Value = 52
$ws.Range("F34").Value = 283
$ws.Range("F35").Value = 47
$ws.Range("F36").Value = 92
$ws.Range("F38").Value = 1233
$ws.Range("F40").Value = 77
$ws.Range("F41").Value = 103
$ws.Range("F43").Value = 2592
$ws.Range("F44").Value = 936
$ws.Range("F45").Value = 317
$ws.Range("F48").Value = 30
$ws.Range("F50").Value = 113
$ws.Range("F51").Value = 113
